$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.315.92"
$ws.Range("E2").Value = '  +1.83%  '

$ws.Range("D3").Value = "'1.647.43"
$ws.Range("E3").Value = '  +0.35%  '

$ws.Range("E4").Value = '  -0.20%  '

$ws.Range("D5").Value = "'217.52"
$ws.Range("E5").Value = '  +0.61%  '

$ws.Range("E6").Value = '  +0.26%  '

$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = '  -0.22%  '

$ws.Range("E8").Value = '  +0.17%  '

$ws.Range("E9").Value = '  +0.23%  '

$ws.Range("D10").Value = "'20.02"
$ws.Range("E10").Value = '  +1.41%  '

$ws.Range("E11").Value = '  +0.11%  '

$ws.Range("E12").Value = '  +0.47%  '

$ws.Range("E13").Value = '  +0.75%  '

$ws.Range("D14").Value = "'1.654.30"
$ws.Range("E14").Value = '  +0.81%  '

$ws.Range("E15").Value = '  -1.96%  '

$ws.Range("E16").Value = '  -0.32%  '

$ws.Range("E17").Value = '  +0.78%  '

$ws.Range("D18").Value = "'26.297.54"
$ws.Range("E18").Value = '  +1.59%  '

$ws.Range("D19").Value = "'1.00"
$ws.Range("E19").Value = '  -0.21%  '

$ws.Range("D20").Value = "'196.89"
$ws.Range("E20").Value = '  +1.99%  '

$ws.Range("D22").Value = "'10.08"
$ws.Range("E22").Value = '  +1.12%  '

$ws.Range("D23").Value = "'6.34"
$ws.Range("E23").Value = '  -0.08%  '

$ws.Range("E24").Value = '  -2.64%  '

$ws.Range("D25").Value = "'143.09"
$ws.Range("E25").Value = '  +0.63%  '

$ws.Range("E26").Value = '  -0.17%  '

$ws.Range("E27").Value = '  +1.67%  '

$ws.Range("E28").Value = '  +0.30%  '

$ws.Range("D29").Value = "'15.69"
$ws.Range("E29").Value = '  +0.84%  '

$ws.Range("D30").Value = "'1.25"
$ws.Range("E30").Value = '  +1.05%  '

$ws.Range("E31").Value = '  +2.89%  '

$ws.Range("D32").Value = "'3.35"
$ws.Range("E32").Value = '  +0.41%  '

$ws.Range("D33").Value = "'3.25"
$ws.Range("E33").Value = '  +0.02%  '

$ws.Range("E34").Value = '  +1.97%  '

$ws.Range("E35").Value = '  +0.99%  '

$ws.Range("E36").Value = '  +0.69%  '

$ws.Range("D37").Value = "'1.138.45"
$ws.Range("E37").Value = '  +0.33%  '

$ws.Range("E38").Value = '  +1.57%  '

$ws.Range("D39").Value = "'2.49"
$ws.Range("E39").Value = '  -1.71%  '

$ws.Range("E40").Value = '  +0.77%  '

$ws.Range("E41").Value = '  -0.24%  '

$ws.Range("D42").Value = "'5.66"
$ws.Range("E42").Value = '  +1.46%  '

$ws.Range("D43").Value = "'100.44"
$ws.Range("E43").Value = '  -0.37%  '

$ws.Range("E44").Value = '  -0.47%  '

$ws.Range("D45").Value = "'1.786.55"
$ws.Range("E45").Value = '  +0.52%  '

$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D46").Value = "'56.50"
$ws.Range("E46").Value = '  +1.87%  '

$ws.Range("B47").Value = 'RenderToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D47").Value = "'1.49"
$ws.Range("E47").Value = '  +3.23%  '

$ws.Range("B48").Value = 'Cronos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D48").Value = "'0.0517"
$ws.Range("E48").Value = '  +2.83%  '

$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").Value = "'7.75"
$ws.Range("E49").Value = '  +3.53%  '

$ws.Range("B50").Value = 'Mantle'
$ws.Range("C50").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D50").Value = "'0.417"
$ws.Range("E50").Value = '  -0.32%  '

$ws.Range("B51").Value = 'Algorand'
$ws.Range("C51").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D51").Value = "'0.0977"
$ws.Range("E51").Value = '  +2.24%  '
